# cryptos.xlsx update — Mon Dec  4 13:29:54 UTC 2023 GitHub Actions refresh.
# Each target cell is forced to Text (NumberFormat "@") before the new value is
# written so numeric-looking strings (e.g. "229.78") stay text like the original
# inlineStr cells, then the style is reset back to Normal so no visible
# formatting/style is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = '@'
    $rng.Value = $val
    $rng.Style = 'Normal'
}

Set-TextValue 'D2' '41.914.93'
Set-TextValue 'E2' '  +5.76%  '
Set-TextValue 'D3' '2.248.89'
Set-TextValue 'E3' '  +3.98%  '
Set-TextValue 'E4' '  +0.21%  '
Set-TextValue 'D5' '229.78'
Set-TextValue 'E5' '  +0.69%  '
Set-TextValue 'D6' '0.628'
Set-TextValue 'E6' '  -0.60%  '
Set-TextValue 'D7' '61.30'
Set-TextValue 'E7' '  -3.37%  '
Set-TextValue 'D9' '0.406'
Set-TextValue 'E9' '  +3.27%  '
Set-TextValue 'D10' '58.58'
Set-TextValue 'E10' '  +0.91%  '
Set-TextValue 'D11' '0.0880'
Set-TextValue 'E11' '  +3.58%  '
Set-TextValue 'E12' '  +0.39%  '
Set-TextValue 'D13' '2.583.26'
Set-TextValue 'E13' '  +4.00%  '
Set-TextValue 'D14' '15.87'
Set-TextValue 'E14' '  -0.69%  '
Set-TextValue 'D15' '21.74'
Set-TextValue 'E15' '  -1.11%  '
Set-TextValue 'E16' '  -0.54%  '
Set-TextValue 'D17' '5.62'
Set-TextValue 'E17' '  +2.45%  '
Set-TextValue 'D18' '2.255.25'
Set-TextValue 'E18' '  +4.17%  '
Set-TextValue 'D19' '41.807.20'
Set-TextValue 'D20' '73.32'
Set-TextValue 'E20' '  +1.95%  '
Set-TextValue 'B21' 'Uniswap'
Set-TextValue 'C21' 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue 'D21' '6.15'
Set-TextValue 'E21' '  +0.48%  '
Set-TextValue 'B22' 'ShibaInu'
Set-TextValue 'C22' 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue 'D22' '0.0₃0898'
Set-TextValue 'E22' '  +6.16%  '
Set-TextValue 'D23' '249.14'
Set-TextValue 'E23' '  +9.17%  '
Set-TextValue 'E24' '  +0.09%  '
Set-TextValue 'D25' '2.38'
Set-TextValue 'E25' '  +0.16%  '
Set-TextValue 'D26' '2.39'
Set-TextValue 'E26' '  -0.33%  '
Set-TextValue 'D27' '9.67'
Set-TextValue 'E27' '  +0.33%  '
Set-TextValue 'D28' '0.142'
Set-TextValue 'E28' '  +3.00%  '
Set-TextValue 'D29' '168.07'
Set-TextValue 'E29' '  -2.38%  '
Set-TextValue 'D30' '20.18'
Set-TextValue 'E30' '  +1.99%  '
Set-TextValue 'E31' '  +2.05%  '
Set-TextValue 'E32' '  +4.72%  '
Set-TextValue 'E33' '  +0.39%  '
Set-TextValue 'D34' '5.16'
Set-TextValue 'E34' '  +9.94%  '
Set-TextValue 'D35' '4.70'
Set-TextValue 'E35' '  +1.97%  '
Set-TextValue 'E36' '  +1.48%  '
Set-TextValue 'D37' '3.81'
Set-TextValue 'E37' '  +5.10%  '
Set-TextValue 'E38' '  -3.26%  '
Set-TextValue 'D39' '2.40'
Set-TextValue 'E39' '  -0.20%  '
Set-TextValue 'D40' '0.000252'
Set-TextValue 'E40' '  +33.40%  '
Set-TextValue 'E41' '  +0.19%  '
Set-TextValue 'D42' '4.94'
Set-TextValue 'E42' '  +3.13%  '
Set-TextValue 'D43' '0.0238'
Set-TextValue 'E43' '  +5.26%  '
Set-TextValue 'D44' '8.77'
Set-TextValue 'E44' '  +13.23%  '
Set-TextValue 'D45' '100.78'
Set-TextValue 'E45' '  -1.28%  '
Set-TextValue 'D46' '0.0982'
Set-TextValue 'E46' '  +6.14%  '
Set-TextValue 'D47' '1.485.81'
Set-TextValue 'E47' '  -1.76%  '
Set-TextValue 'E48' '  -1.76%  '
Set-TextValue 'D49' '16.69'
Set-TextValue 'E49' '  -5.42%  '
Set-TextValue 'E50' '  +0.37%  '
Set-TextValue 'D51' '2.79'
Set-TextValue 'E51' '  -0.46%  '
